$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = '@'
    $r.Value = $value
    $r.Style = 'Normal'
}

Set-TextValue 'D2' '30.309.98'
Set-TextValue 'E2' '  -2.37%  '
Set-TextValue 'D3' '1.932.25'
Set-TextValue 'E3' '  -2.29%  '
Set-TextValue 'D4' '1.002'
Set-TextValue 'E4' '  -0.50%  '
Set-TextValue 'D5' '250.71'
Set-TextValue 'E5' '  -0.59%  '
Set-TextValue 'D6' '0.7115'
Set-TextValue 'E6' '  -2.04%  '
Set-TextValue 'D7' '1.000'
Set-TextValue 'E7' '  -0.44%  '
Set-TextValue 'D8' '0.3296'
Set-TextValue 'E8' '  -1.56%  '
Set-TextValue 'D9' '27.52'
Set-TextValue 'E9' '  -0.09%  '
Set-TextValue 'D10' '0.07239'
Set-TextValue 'E10' '  +2.22%  '
Set-TextValue 'D11' '0.8036'
Set-TextValue 'E11' '  -2.30%  '
Set-TextValue 'D12' '0.08086'
Set-TextValue 'E12' '  +0.07%  '
Set-TextValue 'D13' '1.932.09'
Set-TextValue 'E13' '  -2.46%  '
Set-TextValue 'D14' '5.460'
Set-TextValue 'E14' '  -1.61%  '
Set-TextValue 'D15' '94.37'
Set-TextValue 'E15' '  -4.28%  '
Set-TextValue 'D16' '15.03'
Set-TextValue 'E16' '  -1.20%  '
Set-TextValue 'D17' '30.305.97'
Set-TextValue 'E17' '  -2.36%  '
Set-TextValue 'D18' '252.28'
Set-TextValue 'E18' '  -5.54%  '
Set-TextValue 'D19' '0.000008165'
Set-TextValue 'E19' '  +0.19%  '
Set-TextValue 'D20' '5.782'
Set-TextValue 'E20' '  -4.13%  '
Set-TextValue 'D21' '2.187.59'
Set-TextValue 'E21' '  -2.82%  '
Set-TextValue 'E22' '  -0.51%  '
Set-TextValue 'D23' '1.001'
Set-TextValue 'E23' '  -0.62%  '
Set-TextValue 'D24' '6.963'
Set-TextValue 'E24' '  -0.97%  '
Set-TextValue 'D25' '9.725'
Set-TextValue 'E25' '  -1.79%  '
Set-TextValue 'D26' '165.18'
Set-TextValue 'E26' '  +2.11%  '
Set-TextValue 'D27' '2.347'
Set-TextValue 'E27' '  +0.34%  '
Set-TextValue 'D28' '19.27'
Set-TextValue 'E28' '  -1.67%  '
Set-TextValue 'D29' '0.1285'
Set-TextValue 'E29' '  -2.70%  '
Set-TextValue 'D30' '1.351'
Set-TextValue 'E30' '  -1.51%  '
Set-TextValue 'E31' '  -2.55%  '
Set-TextValue 'D32' '4.414'
Set-TextValue 'E32' '  -3.93%  '
Set-TextValue 'D33' '4.168'
Set-TextValue 'E33' '  -4.91%  '
Set-TextValue 'D34' '0.05188'
Set-TextValue 'E34' '  -1.52%  '
Set-TextValue 'D35' '1.258'
Set-TextValue 'E35' '  -0.92%  '
Set-TextValue 'D36' '0.7447'
Set-TextValue 'E36' '  -3.77%  '
Set-TextValue 'D37' '2.792'
Set-TextValue 'E37' '  +0.28%  '
Set-TextValue 'D38' '0.01963'
Set-TextValue 'E38' '  -1.53%  '
Set-TextValue 'D39' '2.807'
Set-TextValue 'E39' '  -2.33%  '
Set-TextValue 'D40' '78.61'
Set-TextValue 'E40' '  -5.26%  '
Set-TextValue 'D41' '6.413'
Set-TextValue 'E41' '  -4.40%  '
Set-TextValue 'D42' '0.4514'
Set-TextValue 'E42' '  -1.93%  '
Set-TextValue 'D43' '2.013'
Set-TextValue 'E43' '  -3.20%  '
Set-TextValue 'D44' '0.8442'
Set-TextValue 'E44' '  -0.80%  '
Set-TextValue 'E45' '  -0.46%  '
Set-TextValue 'D46' '101.39'
Set-TextValue 'E46' '  -2.90%  '
Set-TextValue 'D47' '9.796'
Set-TextValue 'E47' '  -2.78%  '
Set-TextValue 'D48' '7.427'
Set-TextValue 'E48' '  -2.46%  '
Set-TextValue 'D49' '36.66'
Set-TextValue 'E49' '  -1.09%  '
Set-TextValue 'D50' '0.4161'
Set-TextValue 'E50' '  -2.59%  '
Set-TextValue 'D51' '0.06030'
